$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2-9 down to 3-10)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the NFE11-TP-REG record
$ws.Cells.Item(2, 1).Value = "NFE11-TP-REG"
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 2
$ws.Cells.Item(2, 4).Value = "NUMERO"
$ws.Cells.Item(2, 5).Value = "N"
$ws.Cells.Item(2, 6).Value = ""
